$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "49+49=98",
    "9+64=73",
    "54-7=47",
    "87+8=95",
    "65-56=9",
    "43-5=38",
    "25-16=9",
    "85-17=68",
    "28+67=95",
    "95-28=67",
    "65+6=71",
    "46+17=63",
    "83-64=19",
    "91-62=29",
    "60-26=34",
    "28+37=65",
    "63-55=8",
    "25+27=52",
    "90-68=22",
    "45-38=7",
    "83-79=4",
    "34+17=51",
    "67-39=28",
    "71-42=29",
    "74-16=58",
    "28+59=87",
    "8+39=47",
    "30-25=5",
    "47+4=51",
    "34-6=28",
    "25+59=84",
    "92-34=58",
    "57+29=86",
    "69+2=71",
    "85-26=59",
    "91-17=74",
    "14+39=53",
    "85-38=47",
    "38+17=55",
    "76-38=38",
    "74+8=82",
    "90-85=5",
    "5+59=64",
    "82-15=67",
    "6+66=72",
    "20-11=9",
    "40-39=1",
    "34+47=81",
    "15+28=43",
    "65-47=18",
    "46-17=29",
    "10-9=1",
    "18+76=94",
    "73-16=57",
    "18+36=54",
    "85-48=37",
    "84+7=91",
    "68+16=84",
    "29+62=91",
    "30-9=21",
    "52+29=81",
    "13+59=72",
    "32-8=24",
    "17+78=95",
    "60-35=25",
    "39+26=65",
    "48+44=92",
    "95-26=69",
    "71-53=18",
    "59+6=65",
    "75-9=66",
    "49+29=78",
    "40-9=31",
    "16-8=8",
    "26-8=18",
    "25+68=93",
    "94-77=17",
    "91-55=36",
    "60-13=47",
    "44+38=82",
    "38+49=87",
    "71-47=24",
    "39+58=97",
    "19+53=72",
    "57+8=65",
    "29+33=62",
    "39+24=63",
    "8+35=43",
    "71-47=24",
    "3+79=82",
    "46-18=28",
    "11-9=2",
    "60-47=13",
    "90-29=61",
    "88+8=96",
    "57-18=39",
    "18+5=23",
    "77+6=83",
    "59+13=72",
    "58+18=76"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "cells"
